# Hortaliza, Vega Modelo de Temuco - Pepino dulce: weekly data refresh.
# A new week's record is inserted at row 237 (pushing the existing rows
# 237-358 down to 238-359), and the new row carries this week's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 237; this shifts rows 237:358
# down to 238:359 (formats/values move with them), growing the sheet's
# used range from A1:R358 to A1:R359.
$ws.Rows.Item(237).Insert()

# Populate the newly inserted row 237 with this week's observation.
$ws.Range("A237").Value = 10
$ws.Range("B237").Value = "Vega Modelo de Temuco"
$ws.Range("C237").Value = "La Araucanía"
$ws.Range("D237").Value = 45089
$ws.Range("E237").Value = 9
$ws.Range("F237").Value = 100112043
$ws.Range("G237").Value = "Pepino dulce"
$ws.Range("H237").Value = "Cultivar IV Región"
$ws.Range("I237").Value = "Primera"
$ws.Range("J237").Value = 205
$ws.Range("K237").Value = 16000
$ws.Range("L237").Value = 17000
$ws.Range("M237").Value = 16390
$ws.Range("N237").Value = "`$/bandeja 18 kilos"
$ws.Range("O237").Value = "Provincia de Limarí"
$ws.Range("P237").Value = 911
$ws.Range("Q237").Value = 18
$ws.Range("R237").Value = "Hortaliza"
